# update IMGT result parser
# later need to confirm the germline results with more testing data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: new "function" flag cell + a numeric ID cell -----------------
$ws.Range("F14").Value = "function"
$ws.Range("F14").Font.Color = 255          # reuse red-font style (matches existing s=2)

$ws.Range("H14").Value = 12
$ws.Range("H14").Interior.Color = 65535    # reuse yellow-fill style (matches existing s=1)

# --- Rows 61-68: add "function" marker column (F) + blank flag column (G) -
# These cells all pick up one brand-new (shared) style.
foreach ($r in 61..68) {
    $ws.Range("F$r").Value = "function"
    $ws.Range("F$r").Font.Bold = $true
    $ws.Range("G$r").Font.Bold = $true
}

# --- Row 81: G81 switches from a number to the new text "6 pr 7" ----------
$ws.Range("G81").Value = "6 pr 7"

# --- Row 82: add E82 ("yes") and F82 ("function", new shared style) -------
$ws.Range("E82").Value = "yes"
$ws.Range("F82").Value = "function"
$ws.Range("F82").Font.Bold = $true

# --- Row 99: restyle the existing G99/H99 cells ----------------------------
$ws.Range("G99").Font.Color = 255
$ws.Range("H99").Font.Color = 255

# --- Row 103: restyle F103/H103, add styled blank G103 --------------------
$ws.Range("F103").Font.Color = 255
$ws.Range("G103").Font.Color = 255
$ws.Range("H103").Font.Color = 255

# --- Cosmetic view/page state (best effort) --------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("L81").Select()
